$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.673.44"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.287.82"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.85%  "
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "2.618.51"
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").Value = "2.292.95"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "43.600.61"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000108"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.37%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0884"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  -4.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0354"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -13.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.236"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.32%  "
$ws.Range("E42").Value = "  +16.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  +8.18%  "
$ws.Range("D51").Value = "2.507.19"
$ws.Range("E51").Value = "  +2.21%  "
